# "Automatic update of files."
# Column C ("Förändrad" / Changed-date) for rows 2-10 moves from
# serial date 45224 (2023-10-25) to 45233 (2023-11-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    $old = $cell.Value()

    $isOldDate = $false
    if ($old -is [DateTime]) {
        if ($old.Year -eq 2023 -and $old.Month -eq 10 -and $old.Day -eq 25) {
            $isOldDate = $true
        }
    } elseif ($old -eq 45224) {
        $isOldDate = $true
    }

    if ($isOldDate) {
        $cell.Value = 45233
    }
}
